$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text (non-numeric) cell updates: coin name / link swaps for rows 41-43 ---
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'

# --- Numeric-looking text cell updates (quote-prefixed to keep them stored as text) ---
$ws.Range("D2").Value = "'245.34"
$ws.Range("E2").Value = "'-0.52%"
$ws.Range("G2").Value = "'18"
$ws.Range("E3").Value = "'1.86%"
$ws.Range("G3").Value = "'18"
$ws.Range("D4").Value = "'5.071"
$ws.Range("E4").Value = "'-0.23%"
$ws.Range("G4").Value = "'18"
$ws.Range("D5").Value = "'0.05696"
$ws.Range("E5").Value = "'1.45%"
$ws.Range("G5").Value = "'18"
$ws.Range("D6").Value = "'6.488"
$ws.Range("E6").Value = "'0.20%"
$ws.Range("G6").Value = "'18"
$ws.Range("D7").Value = "'0.8201"
$ws.Range("E7").Value = "'0.80%"
$ws.Range("G7").Value = "'18"
$ws.Range("D8").Value = "'0.8419"
$ws.Range("E8").Value = "'-0.43%"
$ws.Range("G8").Value = "'18"
$ws.Range("D9").Value = "'0.1331"
$ws.Range("E9").Value = "'-0.54%"
$ws.Range("G9").Value = "'18"
$ws.Range("D10").Value = "'0.06916"
$ws.Range("E10").Value = "'-0.62%"
$ws.Range("G10").Value = "'18"
$ws.Range("D11").Value = "'0.02826"
$ws.Range("E11").Value = "'-0.66%"
$ws.Range("G11").Value = "'18"
$ws.Range("D12").Value = "'0.09398"
$ws.Range("E12").Value = "'0.09%"
$ws.Range("G12").Value = "'18"
$ws.Range("D13").Value = "'0.001514"
$ws.Range("E13").Value = "'0.13%"
$ws.Range("G13").Value = "'18"
$ws.Range("D14").Value = "'0.04091"
$ws.Range("E14").Value = "'-12.20%"
$ws.Range("G14").Value = "'18"
$ws.Range("D15").Value = "'0.0006027"
$ws.Range("E15").Value = "'1.10%"
$ws.Range("G15").Value = "'18"
$ws.Range("D16").Value = "'0.006057"
$ws.Range("E16").Value = "'-0.87%"
$ws.Range("G16").Value = "'18"
$ws.Range("D17").Value = "'3.512"
$ws.Range("E17").Value = "'-2.46%"
$ws.Range("G17").Value = "'18"
$ws.Range("D18").Value = "'3.005"
$ws.Range("E18").Value = "'-0.18%"
$ws.Range("G18").Value = "'18"
$ws.Range("D19").Value = "'2.229"
$ws.Range("E19").Value = "'8.47%"
$ws.Range("G19").Value = "'18"
$ws.Range("G20").Value = "'18"
$ws.Range("D21").Value = "'0.03180"
$ws.Range("E21").Value = "'-0.47%"
$ws.Range("G21").Value = "'18"
$ws.Range("D22").Value = "'0.1274"
$ws.Range("E22").Value = "'-0.66%"
$ws.Range("G22").Value = "'18"
$ws.Range("D23").Value = "'3.550"
$ws.Range("E23").Value = "'-5.52%"
$ws.Range("G23").Value = "'18"
$ws.Range("E24").Value = "'1.70%"
$ws.Range("G24").Value = "'18"
$ws.Range("D25").Value = "'0.001217"
$ws.Range("E25").Value = "'-2.65%"
$ws.Range("G25").Value = "'18"
$ws.Range("D26").Value = "'0.003969"
$ws.Range("E26").Value = "'-13.44%"
$ws.Range("G26").Value = "'18"
$ws.Range("D27").Value = "'0.00009789"
$ws.Range("E27").Value = "'1.96%"
$ws.Range("G27").Value = "'18"
$ws.Range("D28").Value = "'0.0001448"
$ws.Range("E28").Value = "'-25.29%"
$ws.Range("G28").Value = "'18"
$ws.Range("G29").Value = "'18"
$ws.Range("G30").Value = "'18"
$ws.Range("G31").Value = "'18"
$ws.Range("G32").Value = "'18"
$ws.Range("G33").Value = "'18"
$ws.Range("G34").Value = "'18"
$ws.Range("G35").Value = "'18"
$ws.Range("G36").Value = "'18"
$ws.Range("G37").Value = "'18"
$ws.Range("G38").Value = "'18"
$ws.Range("G39").Value = "'18"
$ws.Range("D40").Value = "'0.03692"
$ws.Range("E40").Value = "'0.72%"
$ws.Range("G40").Value = "'18"
$ws.Range("D41").Value = "'0.1058"
$ws.Range("E41").Value = "'-22.76%"
$ws.Range("G41").Value = "'18"
$ws.Range("D42").Value = "'0.002381"
$ws.Range("E42").Value = "'-10.49%"
$ws.Range("G42").Value = "'18"
$ws.Range("D43").Value = "'0.003429"
$ws.Range("E43").Value = "'-44.88%"
$ws.Range("G43").Value = "'18"
$ws.Range("D44").Value = "'0.009395"
$ws.Range("E44").Value = "'7.29%"
$ws.Range("G44").Value = "'18"
$ws.Range("E45").Value = "'-1.82%"
$ws.Range("G45").Value = "'18"
$ws.Range("E46").Value = "'-0.07%"
$ws.Range("G46").Value = "'18"
$ws.Range("E47").Value = "'-0.07%"
$ws.Range("G47").Value = "'18"
$ws.Range("D48").Value = "'0.002463"
$ws.Range("E48").Value = "'7.22%"
$ws.Range("G48").Value = "'18"
$ws.Range("E49").Value = "'-0.07%"
$ws.Range("G49").Value = "'18"
$ws.Range("E50").Value = "'-0.07%"
$ws.Range("G50").Value = "'18"
$ws.Range("G51").Value = "'18"
